# Add a new "March-25 (11)" sheet, cloned from "Feb-25 (10)", for the
# next month's (blank) daily purchase report, and make it the active tab.

$wb = $excel.ActiveWorkbook

# Clone the Feb-25 (10) sheet (keeps all formatting/formulas) and place
# the copy right after it.
$template = $wb.Worksheets.Item("Feb-25 (10)")
$template.Copy($null, $template)

# The newly created copy is now the last sheet.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "March-25 (11)"

# It's a brand-new, not-yet-filled-in report: wipe out last month's
# purchase figures (columns C..R for the 31 day rows), leaving the
# totals/formulas in row 34+ intact (they'll just recompute to 0).
$newSheet.Range("C3:R33").ClearContents()

# Re-date the 31 day rows for March 2025 (1-Mar-25 .. 31-Mar-25).
for ($i = 0; $i -lt 31; $i++) {
    $row = 3 + $i
    $newSheet.Cells.Item($row, 2).Value = 45352 + $i
}

# Match the author's saved selection/view state on the new sheet.
$newSheet.Range("A6").Select()
